$d = $word.ActiveDocument
Write-Host "Paragraphs:" $d.Paragraphs.Count
